$d = $word.ActiveDocument

$replacements = @(
    @{old="333÷9="; new="419÷2="},
    @{old="379÷8="; new="574÷5="},
    @{old="997÷3="; new="783÷9="},
    @{old="967÷5="; new="878÷6="},
    @{old="726÷5="; new="637÷2="},
    @{old="850÷5="; new="677÷4="},
    @{old="647÷5="; new="613÷3="},
    @{old="445÷3="; new="500÷6="},
    @{old="751÷2="; new="167÷9="},
    @{old="941÷8="; new="337÷9="},
    @{old="455÷2="; new="359÷3="},
    @{old="493÷8="; new="503÷5="},
    @{old="562÷9="; new="442÷7="},
    @{old="379÷9="; new="252÷7="},
    @{old="866÷2="; new="784÷3="},
    @{old="453÷3="; new="192÷6="},
    @{old="520÷4="; new="117÷9="},
    @{old="986÷4="; new="702÷9="},
    @{old="876÷5="; new="257÷2="},
    @{old="177÷9="; new="856÷2="},
    @{old="595÷3="; new="654÷7="},
    @{old="228÷4="; new="266÷8="},
    @{old="405÷6="; new="480÷3="},
    @{old="753÷9="; new="320÷6="},
    @{old="924÷9="; new="257÷8="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
